## manual test for index 2020-10-28 23:11:09
## Appends a duplicate of the last row (A50 -> A51) to the single-column
## "indexed" list sheet, keeping the same shared-string value/style, then
## moves the selection/scroll position to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last populated row currently is 50; the new row goes right after it.
$lastRow = 50
$newRow  = $lastRow + 1

$src = $ws.Cells.Item($lastRow, 1)
$dst = $ws.Cells.Item($newRow, 1)

# Copy (value + style) so the new cell reuses the existing shared-string
# entry and the Arial cell style instead of minting new ones.
$src.Copy($dst)

# Move the selection to the freshly written cell and scroll the window
# down so row 22 is at the top, matching where Excel would land after
# typing into the next free row of a long list.
$dst.Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
